$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.732.81'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '2.039.35'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'227.26"
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('E6').Value = '  -1.12%  '
$ws.Range('D7').Value = "'59.60"
$ws.Range('E7').Value = '  -0.83%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = "'0.377"
$ws.Range('E9').Value = '  -2.46%  '
$ws.Range('D10').Value = "'0.0840"
$ws.Range('E10').Value = '  +2.88%  '
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('D12').Value = '2.340.46'
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('E13').Value = '  -1.26%  '
$ws.Range('E14').Value = '  -0.54%  '
$ws.Range('D15').Value = "'5.46"
$ws.Range('E15').Value = '  +4.39%  '
$ws.Range('D16').Value = "'0.772"
$ws.Range('E16').Value = '  +2.27%  '
$ws.Range('D17').Value = '2.042.47'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').Value = '37.646.89'
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('D19').Value = "'69.45"
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('D20').Value = "'5.91"
$ws.Range('E20').Value = '  -1.68%  '
$ws.Range('D21').Value = '0.0₃0824'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').Value = "'223.55"
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('E25').Value = '  +2.36%  '
$ws.Range('D26').Value = "'168.24"
$ws.Range('E26').Value = '  +2.06%  '
$ws.Range('D27').Value = "'9.38"
$ws.Range('E27').Value = '  +1.41%  '
$ws.Range('D28').Value = "'0.128"
$ws.Range('E28').Value = '  -1.13%  '
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('E31').Value = '  -0.56%  '
$ws.Range('E32').Value = '  +8.18%  '
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('D34').Value = "'0.0605"
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('D35').Value = "'4.49"
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('D36').Value = "'6.48"
$ws.Range('E36').Value = '  +1.74%  '
$ws.Range('D37').Value = "'2.33"
$ws.Range('E37').Value = '  +3.30%  '
$ws.Range('D38').Value = "'3.41"
$ws.Range('E38').Value = '  +5.00%  '
$ws.Range('D39').Value = "'0.999"
$ws.Range('E39').Value = '  -0.17%  '
$ws.Range('D40').Value = "'18.07"
$ws.Range('E40').Value = '  +8.56%  '
$ws.Range('D41').Value = '1.528.53'
$ws.Range('E41').Value = '  -0.55%  '
$ws.Range('D42').Value = "'97.29"
$ws.Range('E42').Value = '  +0.45%  '
$ws.Range('E43').Value = '  -0.97%  '
$ws.Range('E44').Value = '  +1.30%  '
$ws.Range('D45').Value = "'4.24"
$ws.Range('E45').Value = '  +6.91%  '
$ws.Range('E46').Value = '  -1.64%  '
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').Value = "'1.00"
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('E50').Value = '  -1.76%  '
$ws.Range('D51').Value = '2.228.80'
$ws.Range('E51').Value = '  +0.42%  '
